# Ryan Cromar email 11/2/2021 1:24 PM
#
# 1) Question 10 answer: "A Q-Q Plot" -> "A histogram"
# 2) Part (b) of the non-certified/certified cars answer: reworded to
#    talk about the sampling distribution of the sample mean, and to
#    reference a histogram (instead of a Q-Q Plot) for checking normality.

$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "10 . A Q-Q Plot will help us know if the data are normally distributed.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "10 . A histogram will help us know if the data are normally distributed.",
    2) | Out-Null

$d.Content.Find.Execute(
    "b. The sample size for non-certified cars is large, so we can assume that it is normally distributed. With a sample size of n = 24 for certified cars, it is likely large enough to assume a normal distribution, but a Q-Q Plot will help us know more confidently.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "b. The sample size for non-certified cars is large, so we can assume that the sampling distribution of the sample mean is normally distributed. With a sample size of n = 24 for certified cars, it is likely large enough to assume a normal distribution for the mean, but a histogram of the data will help us know if the data is normally distributed. If it is, then the distribution of the sample mean will also be normally distributed.",
    2) | Out-Null
